$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Corrijo rnf": the alternative in the Java decision row was just the bare
# word "Java"; reword it to a fuller description.
$ws.Range("B2").Value = "Lenguaje de programacion Java"

# "Agrego cotas": row 3 (the "cotas"/categories decision) no longer needs as
# much vertical space once the table re-flowed, so its height shrinks.
$ws.Rows.Item(3).RowHeight = 26.25

# Leave the selection where the author's cursor ended up on save.
$ws.Range("D5").Select()
